$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: new product "пылесос" (vacuum cleaner)
$ws.Range("A6").Value = "пылесос"
$ws.Range("B6").Value = "https://uz.ozon.com/product/pylesos-vertikalnyy-besprovodnoy-futula-v8-belyy-vysokaya-moshchnost-v-9500pa-turbo-shchetka-s-1348678583/?_bctx=CAYQ-98y&at=x6tPEj4LZhgAGLP1TY887Q7Cn9w2xOc502K98UnkLQj6"
$ws.Range("C6").Value = '//*[@id="layoutPage"]/div[1]/div[3]/div[3]/div[2]/div/div/div[1]/div[2]/div/div[1]/div/div/div[1]/div/div[1]/div[1]/span[1]'
$ws.Rows.Item(6).RowHeight = 210

# Row 7: new product "сумка" (bag)
$ws.Range("A7").Value = "сумка"
$ws.Range("B7").Value = "https://aliexpress.ru/item/1005008688882796.html?spm=a2g2w.home.3.1.139f5586h1RuUJ&mixer_rcmd_bucket_id=aerabtestalgoRecommendAbV2_testRankingFairPriceMerged&ru_algo_pv_id=5f2048-ade163-e7d317-9725ba-1744909200&scenario=aerMediaKitSegments&spmC=homepage_main_shelf_pc&spmD=3&traffic_source=recommendation&type_rcmd=core&sku_id=12000046251455202"
$ws.Range("C7").Value = '//*[@id="__aer_root__"]/div/div[8]/div[2]/div[4]/div/div[1]/div[1]/div/div[1]/div/div/div[2]'
$ws.Rows.Item(7).RowHeight = 165

# Scroll the view down and move the active selection to E7, matching the
# author's final cursor position after adding the two new rows.
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("E7").Select()
